$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at 112-114 (pushes the existing 112-114 rows down to 115-117)
$ws.Range("A112:A114").EntireRow.Insert()

# Grow the table (ListObject) to cover the 3 new rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E117"))

# Row 112: MulticastKitDeleted
$ws.Cells.Item(112,1).Value = "MulticastKitDeleted"
$ws.Cells.Item(112,2).Value = 3006
$ws.Cells.Item(112,3).Value = "KitSync.NetCalls"
$ws.Cells.Item(112,4).Value = "FROM_EITHER"
$ws.Cells.Item(112,5).Value = "int pk"

# Row 113: MulticastKitAccessChanged
$ws.Cells.Item(113,1).Value = "MulticastKitAccessChanged"
$ws.Cells.Item(113,2).Value = 3007
$ws.Cells.Item(113,3).Value = "KitSync.NetCalls"
$ws.Cells.Item(113,4).Value = "FROM_EITHER"
$ws.Cells.Item(113,5).Value = "ulong steamId"

# Row 114: MulticastKitUpdated
$ws.Cells.Item(114,1).Value = "MulticastKitUpdated"
$ws.Cells.Item(114,2).Value = 3008
$ws.Cells.Item(114,3).Value = "KitSync.NetCalls"
$ws.Cells.Item(114,4).Value = "FROM_EITHER"
$ws.Cells.Item(114,5).Value = "int pk"

# Update the view state to match the committed workbook
$ws.Application.ActiveWindow.ScrollRow = 85
$ws.Range("E106").Select()
